# Fruta / hortaliza, semanal
# Insert two new weekly price rows at the top of the Frutilla data block
# (rows 619-649 shift down to 621-651), and populate the two new rows
# (619 and 620) with the latest observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 619; existing rows 619:649 shift to 621:651
$ws.Rows("619:620").Insert()

# --- New row 619 ---
$ws.Cells.Item(619, 1).Value = 7
$ws.Cells.Item(619, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(619, 3).Value = "Ñuble"
$ws.Cells.Item(619, 4).Value = 45267
$ws.Cells.Item(619, 5).Value = 16
$ws.Cells.Item(619, 6).Value = "Fruta"
$ws.Cells.Item(619, 7).Value = 100101
$ws.Cells.Item(619, 8).Value = "Berries"
$ws.Cells.Item(619, 9).Value = 100112025
$ws.Cells.Item(619, 10).Value = "Frutilla"
$ws.Cells.Item(619, 11).Value = "Sin especificar"
$ws.Cells.Item(619, 12).Value = "Especial"
$ws.Cells.Item(619, 13).Value = 200
$ws.Cells.Item(619, 14).Value = 12000
$ws.Cells.Item(619, 15).Value = 12000
$ws.Cells.Item(619, 16).Value = 12000
$ws.Cells.Item(619, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(619, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(619, 19).Value = 1714
$ws.Cells.Item(619, 20).Value = 7

# --- New row 620 ---
$ws.Cells.Item(620, 1).Value = 7
$ws.Cells.Item(620, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(620, 3).Value = "Ñuble"
$ws.Cells.Item(620, 4).Value = 45267
$ws.Cells.Item(620, 5).Value = 16
$ws.Cells.Item(620, 6).Value = "Fruta"
$ws.Cells.Item(620, 7).Value = 100101
$ws.Cells.Item(620, 8).Value = "Berries"
$ws.Cells.Item(620, 9).Value = 100112025
$ws.Cells.Item(620, 10).Value = "Frutilla"
$ws.Cells.Item(620, 11).Value = "Sin especificar"
$ws.Cells.Item(620, 12).Value = "Primera"
$ws.Cells.Item(620, 13).Value = 150
$ws.Cells.Item(620, 14).Value = 10000
$ws.Cells.Item(620, 15).Value = 10000
$ws.Cells.Item(620, 16).Value = 10000
$ws.Cells.Item(620, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(620, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(620, 19).Value = 1429
$ws.Cells.Item(620, 20).Value = 7
